$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "Molnar Santoro 2014" row content first (lands on row 8, below the
#     existing "Gaspar et al 2011" row) ---
$ws.Range("A8").Value = "Molnar Santoro 2014"
$ws.Range("B8").Value = "optimal monpol under adaptive EE learning involves a new intertemporal tradeoff in which CB foregoes short-run stabilization in order to facilitate learning"

# --- Insert a new row above "Gaspar et al 2011" (now row 8, shifting it to row 9)
#     for the "Gaspar et al 2006" entry ---
$ws.Range("A7").EntireRow.Insert()
$ws.Range("A7").Value = "Gaspar et al 2006"
$ws.Range("B7").Value = "CB reacts more persistently to cost-push shocks"

# --- Fill in the missing 5-word summary for the existing "Gaspar et al 2011" row ---
$ws.Range("B8").Value = "Hbook chapter"

# --- Append two more rows of new literature ---
$ws.Range("A10").Value = "Ferrero 2007"
$ws.Range("B10").Value = "more aggressive monpol increases the speed of convergence, but high speed is not always desirable"

$ws.Range("A11").Value = "Eusepi et al, 2018, Limits"
$ws.Range("B11").Value = "due to the gain, monetary policy faces limits"

# --- Widen column B to fit the new, longer summaries ---
$ws.Columns.Item(2).ColumnWidth = 133.83072916666666

# --- Leave the selection where the author's cursor ended up after typing ---
$ws.Range("B12").Select()
